# Adds the newest Mega-Sena draws (concursos 2918-2933, rows 375-390) to the
# "MEGA SENA" sheet, keeps the sheet's "most-recent-draws" highlight on the
# new last 10 rows, and updates the view/selection state - matching the
# author's "ajustando o modal dos planos das paginas..." commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MEGA SENA")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1) New draw data (Concurso, Bola1..Bola6) for rows 375-390
# ---------------------------------------------------------------------------
$newDraws = @(
    @(2918, 11, 27, 31, 41, 48, 54),
    @(2919,  3, 26, 28, 37, 42, 53),
    @(2920,  8, 12, 16, 19, 31, 58),
    @(2921,  9, 12, 14, 16, 26, 36),
    @(2922,  4, 23, 30, 39, 40, 41),
    @(2923, 18, 27, 32, 39, 55, 56),
    @(2924, 10, 19, 30, 40, 48, 54),
    @(2925,  7,  9, 12, 13, 24, 27),
    @(2926,  3,  4, 14, 35, 45, 49),
    @(2927, 11, 27, 34, 55, 56, 58),
    @(2928, 14, 24, 29, 32, 46, 48),
    @(2929,  3,  7,  8, 34, 35, 51),
    @(2930,  1, 11, 13, 14, 36, 45),
    @(2931,  4, 19, 23, 36, 47, 52),
    @(2932,  4, 13, 25, 36, 40, 53),
    @(2933,  1, 18, 22, 42, 48, 50)
)

$startRow = 375
for ($i = 0; $i -lt $newDraws.Count; $i++) {
    $row = $startRow + $i
    $values = $newDraws[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# ---------------------------------------------------------------------------
# 2) The "latest draws" highlight always sits on the bottom 10 rows of the
#    table. It now has to move off the old rows (372-380, plain again) and
#    onto the new bottom 10 rows (381-390).
# ---------------------------------------------------------------------------
$ws.Range("A372:G380").Style = "Normal"
$ws.Range("A381:G390").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Selection / scroll position ends on the newly entered block.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 359
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B381:G390").Select()

# ---------------------------------------------------------------------------
# 4) Workbook window geometry (matches the saved bookViews/workbookView).
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = 75285
$win.Top = 675
$win.Width = 18255
$win.Height = 17505
